$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right before the existing "2022-Q3"
#    sheet (mirrors Worksheets.Add(Before:=...) in real Excel COM).
# ------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($oldQ3)
$newSheet.Name = "2022-Q4"

# Re-fetch a fresh reference to the "2022-Q3" sheet: the handle obtained
# before Worksheets.Add() can go stale once the sheet collection has been
# mutated, which silently turns a later Copy() into a no-op.
$oldQ3 = $wb.Worksheets.Item("2022-Q3")

# Copy the header row + formatting from the "2022-Q3" sheet (same table
# layout / styles as every quarter sheet) so the new sheet picks up the
# correct column headers, bold/border header style (s=2) and the bold
# index-column style used throughout the workbook.
$oldQ3.Range("A1:H4").Copy($newSheet.Range("A1"))
# The source range has 3 data rows (4 total); the new quarter only has 2
# data rows (3 total), so drop the extra trailing row that Copy brought in.
$newSheet.Rows.Item(4).Delete()

# Re-fetch again after the structural Delete() for the same reason as above.
$newSheet = $wb.Worksheets.Item("2022-Q4")

# The copied range included the (empty) A1 cell, which the source sheet
# never actually populates - drop the leftover placeholder so A1 is blank
# again, matching the other quarter sheets.
$newSheet.Range("A1").ClearContents()

# ------------------------------------------------------------------
# 2. Populate the new "2022-Q4" sheet with the fund holding data.
# ------------------------------------------------------------------
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'014133"
$newSheet.Range("C2").Value = "工银中证500六个月持有指数增强A"
$newSheet.Range("D2").Value = "'1.59"
$newSheet.Range("E2").Value = "'94.29"
$newSheet.Range("F2").Value = "'1.28"
$newSheet.Range("G2").Value = "'0.0204"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'014134"
$newSheet.Range("C3").Value = "工银中证500六个月持有指数增强C"
$newSheet.Range("D3").Value = "'0.88"
$newSheet.Range("E3").Value = "'94.29"
$newSheet.Range("F3").Value = "'1.28"
$newSheet.Range("G3").Value = "'0.0113"
$newSheet.Range("H3").Value = 9

# ------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: add a 2022-Q4 row at the top of
#    the data and push the previously existing rows down by one.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Extend the table by one row, copying the formatting (incl. the bold
# index-column style) of the last existing row down to the new last row.
$summary.Range("A6:D6").Copy($summary.Range("A7"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.03

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 0.1

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q4"
$summary.Range("C4").Value = 5
$summary.Range("D4").Value = 0.23

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.06

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q1"
$summary.Range("C6").Value = 4
$summary.Range("D6").Value = 0.04

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2020-Q4"
$summary.Range("C7").Value = 2
$summary.Range("D7").Value = 0.14
